$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# --- KPI cells: B2 and B3 become blank, currency-formatted ("$"#,##0.00) cells ---
# B3 previously held =SUM(Lotes!H3,Lotes!H100); that formula goes away and the
# cell becomes an (as yet empty) KPI value cell with the currency format.
$ws.Range("B3").ClearContents()
$ws.Range("B3").NumberFormat = """$""#,##0.00"
$ws.Range("B2").NumberFormat = """$""#,##0.00"

# --- B4: new KPI value cell formatted as 0.00 ---
# Apply the 0.00 format to A1 first (mutates its existing style in place),
# copy that resulting format onto B4, then restore A1 back to its default
# (unstyled) appearance - this mirrors how the style got reassigned from
# A1 onto B4 while keeping the same underlying style record.
$ws.Range("A1").NumberFormat = "0.00"
$ws.Range("B4").NumberFormat = $ws.Range("A1").NumberFormat
$ws.Range("A1").ClearFormats()

# --- Move the selection to B5 (new KPI input row) ---
$ws.Range("B5").Select()

# --- Reposition/resize the "Saida por Produtos" chart ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 278.181640625
$co.Top = 8.25
$co.Width = 524.4375
$co.Height = 216
